$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Best-Worst Correl"

# --- Pre-format the previously-empty "Key Biological Features" (N) and
# --- "Potential Reason for Poor Correlation" (Q) columns for the
# --- best/worst-correlation gene block (rows 6-10) so they match the
# --- centered / wrapped style already used by the neighbouring cells.
$ws.Range("P6").Copy()
$ws.Range("N6:N10").PasteSpecial(-4122)
$ws.Range("P6").Copy()
$ws.Range("Q6:Q10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Now fill in the text, in the same order the table was originally authored. ---
$ws.Range("N6").Value  = "Transcriptional regulation of hematopoiesis, T-helper cell differentiation control, histone methylation (H3-K9/H3-K4) regulation"
$ws.Range("Q6").Value  = "Consistent expression in immune cell lineages, regulated epigenetic modification patterns, stable transcriptional activity across developmental stages"

# --- Footnote about table provenance ---
$ws.Range("C12").Value = "I made this table from the separate gene ontology tables of ENRICHR API"

$ws.Range("N7").Value  = "RNA polymerase II CTD phosphorylation, cyclin-dependent protein kinase activity, protein modification"
$ws.Range("Q7").Value  = "Stable expression pattern linked to fundamental transcriptional machinery, essential function in RNA polymerase II regulation"

$ws.Range("N8").Value  = "DNA damage response with p53 signaling, cell cycle regulation, transcriptional activation/repression, apoptotic process"
$ws.Range("Q8").Value  = "Stable expression linked to fundamental cellular processes, consistent role in cell cycle control"

$ws.Range("N9").Value  = "beta/gamma crystallin structural protein, possible transporter function"
$ws.Range("Q9").Value  = "Broad expression pattern across tissues, stable structural protein role, involvement in fundamental cellular transport processes"

$ws.Range("N10").Value = "miRNA catabolic process regulation,  pre-miRNA processing, RNA destabilization, negative regulation of miRNA-mediated gene silencing"
$ws.Range("Q10").Value = "Consistent role in miRNA metabolism, participation in fundamental RNA processing pathways"

# --- Widen N and Q to fit the new text ---
$ws.Columns.Item(14).ColumnWidth = 18.666666666666668
$ws.Columns.Item(17).ColumnWidth = 20.498697916666668

# --- Grow the data rows so the wrapped text is fully visible ---
$ws.Rows.Item(6).RowHeight = 187
$ws.Rows.Item(7).RowHeight = 136
$ws.Rows.Item(8).RowHeight = 119
$ws.Rows.Item(9).RowHeight = 102
$ws.Rows.Item(10).RowHeight = 136

# --- Match the author's final selection ---
$ws.Range("P11").Select()
